$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column B updates ---
$ws.Range("B2").Value = "NSE:ASHIANA"
$ws.Range("B3").Value = "NSE:DPWIRES"
$ws.Range("B4").Value = "NSE:FDC"
$ws.Range("B5").Value = "NSE:GLAND"
$ws.Range("B6").Value = "NSE:INGERRAND"
$ws.Range("B7").Value = "NSE:KOTHARIPRO"
$ws.Range("B8").Value = "NSE:PGIL"
$ws.Range("B9").Value = "NSE:RAMRAT"

# --- Clear column B for rows 10-26 ---
$ws.Range("B10").Value = ""
$ws.Range("B11").Value = ""
$ws.Range("B12").Value = ""
$ws.Range("B13").Value = ""
$ws.Range("B14").Value = ""
$ws.Range("B15").Value = ""
$ws.Range("B16").Value = ""
$ws.Range("B17").Value = ""
$ws.Range("B18").Value = ""
$ws.Range("B19").Value = ""
$ws.Range("B20").Value = ""
$ws.Range("B21").Value = ""
$ws.Range("B22").Value = ""
$ws.Range("B23").Value = ""
$ws.Range("B24").Value = ""
$ws.Range("B25").Value = ""
$ws.Range("B26").Value = ""

# --- Column C updates (row 2 unchanged) ---
$ws.Range("C3").Value = "NSE:3PLAND"
$ws.Range("C4").Value = "NSE:ALPA"
$ws.Range("C5").Value = "NSE:ASMS"
$ws.Range("C6").Value = "NSE:BBL"
$ws.Range("C7").Value = "NSE:CGCL"
$ws.Range("C8").Value = "NSE:CREDITACC"
$ws.Range("C9").Value = "NSE:DEN"
$ws.Range("C10").Value = "NSE:DPABHUSHAN"
$ws.Range("C11").Value = "NSE:EXCELINDUS"
$ws.Range("C12").Value = "NSE:FUSION"
$ws.Range("C13").Value = "NSE:GATEWAY"
$ws.Range("C14").Value = "NSE:GESHIP"
$ws.Range("C15").Value = "NSE:GRASIM"
$ws.Range("C16").Value = "NSE:HERITGFOOD"
$ws.Range("C17").Value = "NSE:ICEMAKE"
$ws.Range("C18").Value = "NSE:JYOTISTRUC"
$ws.Range("C19").Value = "NSE:KALAMANDIR"
$ws.Range("C20").Value = "NSE:LAMBODHARA"
$ws.Range("C21").Value = "NSE:MAPMYINDIA"
$ws.Range("C22").Value = "NSE:MAXIND"
$ws.Range("C23").Value = "NSE:NAVKARCORP"
$ws.Range("C24").Value = "NSE:PFOCUS"
$ws.Range("C25").Value = "NSE:PGHH"
$ws.Range("C26").Value = "NSE:PRIMESECU"

# --- E2 cleared, F2 updated ---
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = "NSE:HINDALCO"

# --- E3 cleared ---
$ws.Range("E3").Value = ""

# --- Delete row 27 (shifts rows up, updates dimension) ---
$ws.Rows(27).Delete()

Write-Output "Edit applied."
